$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 2nd May Data Refresh: update regcntr_id values per refreshed master data
$ws.Range("A3").Value = 10003
$ws.Range("A25").Value = 10003

# Restore/scroll selection state left by the editor after the refresh
$ws.Range("A34:XFD1048576").Select()
